$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{addr='D2'; val='66.239.02'},
    @{addr='E2'; val='  -0.19%  '},
    @{addr='D3'; val='3.565.07'},
    @{addr='E3'; val='  +0.48%  '},
    @{addr='E4'; val='  -0.02%  '},
    @{addr='D5'; val='606.20'},
    @{addr='E5'; val='  -0.39%  '},
    @{addr='D6'; val='144.53'},
    @{addr='E6'; val='  -0.37%  '},
    @{addr='D7'; val='3.563.42'},
    @{addr='E7'; val='  +0.37%  '},
    @{addr='E8'; val='  +0.11%  '},
    @{addr='E9'; val='  +2.26%  '},
    @{addr='E10'; val='  -0.57%  '},
    @{addr='D11'; val='7.83'},
    @{addr='E11'; val='  -2.91%  '},
    @{addr='D12'; val='0.414'},
    @{addr='E12'; val='  -0.33%  '},
    @{addr='D13'; val='4.169.72'},
    @{addr='E13'; val='  +0.56%  '},
    @{addr='E14'; val='  -1.09%  '},
    @{addr='D15'; val='30.37'},
    @{addr='E15'; val='  -0.67%  '},
    @{addr='D16'; val='3.557.78'},
    @{addr='E16'; val='  +0.39%  '},
    @{addr='D17'; val='66.292.62'},
    @{addr='E17'; val='  -0.18%  '},
    @{addr='E18'; val='  -0.63%  '},
    @{addr='D19'; val='11.44'},
    @{addr='E19'; val='  +4.69%  '},
    @{addr='D20'; val='6.22'},
    @{addr='E20'; val='  -0.29%  '},
    @{addr='D21'; val='14.81'},
    @{addr='E21'; val='  -1.44%  '},
    @{addr='D22'; val='430.59'},
    @{addr='E22'; val='  +0.63%  '},
    @{addr='D23'; val='0.614'},
    @{addr='E23'; val='  +1.59%  '},
    @{addr='D24'; val='79.45'},
    @{addr='E24'; val='  +1.13%  '},
    @{addr='D25'; val='3.708.12'},
    @{addr='E25'; val='  +0.57%  '},
    @{addr='E26'; val='  -0.03%  '},
    @{addr='E27'; val='  -3.09%  '},
    @{addr='B28'; val='InternetComputer(DFINITY)'},
    @{addr='C28'; val='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'},
    @{addr='D28'; val='9.20'},
    @{addr='E28'; val='  -1.46%  '},
    @{addr='B29'; val='PancakeSwap'},
    @{addr='C29'; val='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'},
    @{addr='D29'; val='2.50'},
    @{addr='E29'; val='  +0.59%  '},
    @{addr='D30'; val='7.94'},
    @{addr='E30'; val='  -1.87%  '},
    @{addr='D31'; val='1.00'},
    @{addr='E31'; val='  +0.04%  '},
    @{addr='D32'; val='3.559.39'},
    @{addr='E32'; val='  +0.83%  '},
    @{addr='D33'; val='25.44'},
    @{addr='E33'; val='  -0.21%  '},
    @{addr='E34'; val='  -2.88%  '},
    @{addr='E35'; val='  -8.28%  '},
    @{addr='D36'; val='7.88'},
    @{addr='E36'; val='  +0.70%  '},
    @{addr='E38'; val='  -1.72%  '},
    @{addr='E39'; val='  -0.58%  '},
    @{addr='D40'; val='175.04'},
    @{addr='E40'; val='  +2.53%  '},
    @{addr='E41'; val='  -1.40%  '},
    @{addr='D42'; val='5.20'},
    @{addr='E42'; val='  -0.25%  '},
    @{addr='E43'; val='  -0.74%  '},
    @{addr='E44'; val='  +1.79%  '},
    @{addr='D45'; val='45.99'},
    @{addr='E45'; val='  +1.03%  '},
    @{addr='E46'; val='  -0.01%  '},
    @{addr='D47'; val='1.19'},
    @{addr='E47'; val='  -2.13%  '},
    @{addr='D48'; val='2.46'},
    @{addr='E48'; val='  +1.22%  '},
    @{addr='D49'; val='24.93'},
    @{addr='E49'; val='  -4.59%  '},
    @{addr='D50'; val='7.15'},
    @{addr='E50'; val='  -0.90%  '},
    @{addr='D51'; val='23.51'},
    @{addr='E51'; val='  +3.79%  '}
)

foreach ($u in $updates) {
    $ws.Range($u.addr).Value = "'" + $u.val
    $ws.Range($u.addr).Style = "Normal"
}
